# "speed ups to runtime system"
#
# Adds a new benchmark run (JRuby row, cols G:J) to the "fasta" sheet,
# recomputes the comparison figure previously held by the Objeck row,
# converts the plain per-row AVERAGE formulas into a single shared
# formula on each sheet (what Excel does when you fill E2 across
# E2:E11 instead of typing the same formula into every cell), and
# leaves the "fasta" sheet as the active / selected sheet & cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # fannkuch-redux
$ws2 = $wb.Worksheets.Item(2)   # fasta

# ---------------------------------------------------------------
# 1. fannkuch-redux: re-enter E2:E11 as one fill (-> shared formula)
# ---------------------------------------------------------------
$ws1.Range("E2:E11").Formula = "=AVERAGE(B2:D2)"

# Selection on fannkuch-redux moves from B20 to H13:K13, and it stops
# being the tab that's selected/active (fasta becomes active below).
$ws1.Activate()
$ws1.Range("H13:K13").Select()

# ---------------------------------------------------------------
# 2. fasta: re-enter E2:E11 as one fill (-> shared formula), same as above
# ---------------------------------------------------------------
$ws2.Range("E2:E11").Formula = "=AVERAGE(B2:D2)"

# New timing run recorded for the JRuby row (row 3): three raw samples
# plus their average, and a derived "speed up vs. Objeck" percentage
# stashed on row 4.
$ws2.Range("G3").Formula = "=SUM(0*60+53.251)"
$ws2.Range("H3").Formula = "=SUM(0*60+53.068)"
$ws2.Range("I3").Formula = "=SUM(0*60+53.282)"
$ws2.Range("J3").Formula = "=AVERAGE(G3:I3)"
$ws2.Range("J4").Formula = "=1-J3/E7"

# The sorted/pasted summary table (A17:B26) caches values by hand;
# the Objeck row's cached number is refreshed to the new JRuby average.
$ws2.Range("B22").Value = 53.200333333333333

# fasta becomes the active sheet, with C20 selected (was B44).
$ws2.Activate()
$ws2.Range("C20").Select()
